$d = $word.ActiveDocument

# Fix typo: "iearnt" -> "learnt"
$d.Content.Find.Execute("iearnt", $true, $false, $false, $false, $false, $true, 1, $false, "learnt", 2)

# Fix typo: "fo feed" -> "to feed"
$d.Content.Find.Execute("likes fo feed", $true, $false, $false, $false, $false, $true, 1, $false, "likes to feed", 2)

# Fix typo: "the Information above" -> "the information above" (case-sensitive match)
$d.Content.Find.Execute("the Information above", $true, $false, $false, $false, $false, $true, 1, $false, "the information above", 2)

# Remove the paragraph that only contains the inline picture (Picture 1 / 41.jpg)
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0 -and $p.Range.Text.Trim() -eq "") {
        $p.Range.Delete()
    }
}
